# Add 9 new product-category rows (Name / Code) to the bottom of the table,
# matching the new sharedStrings/table/sheet content from the diff, and
# move the visible selection/scroll position to reflect where the user was
# working when they finished (around row 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: Name (column A) / Code (column B)
$newRows = @(
    @("Surgical Kits",            "PROD-051"),
    @("Diagnostic Equipment",     "PROD-052"),
    @("Fitness Equipment",        "PROD-053"),
    @("Diagnostic Accessories",   "PROD-054"),
    @("Patient Care Equipment",   "PROD-055"),
    @("Hospital Accessories",     "PROD-056"),
    @("Physiotherapy Equipment",  "PROD-057"),
    @("Dialysis Equipment",       "PROD-058"),
    @("Emergency Equipment",      "PROD-059")
)

$startRow = 52
$endRow = $startRow + $newRows.Count - 1

# Grow the table ("Table") so the new rows become part of it, matching the
# widened ref="A1:B60" in xl/tables/table1.xml.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B" + $endRow))

# Column B of the existing data rows uses the wrap-text / vertical-center
# style (s="2"). Copy that formatting down onto the new B cells before
# writing values so the new rows keep the same look as the rest of the
# table.
$ws.Range("B51").Copy()
$ws.Range("B" + $startRow + ":B" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write all of the new Name values first, then all of the new Code values,
# so new shared-string entries are appended in the same order the
# original authoring session produced them in (names 102-110, then codes
# 111-119) instead of interleaved per-row.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Reflect the author's final cursor position / scroll: selection on C50,
# scrolled so row 34 is at the top.
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
[void]$ws.Range("C50").Select()
